$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.605.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.604.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.23%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.27"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.06%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.84"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.44%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.38%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.15"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.081.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -8.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.455.03"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.616.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.25%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.55%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.76%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.56%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "588.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.16%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.71%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.160"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.53%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.72"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.83%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.404"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.48%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.65"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.49"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.28"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "156.09"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.14%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.11"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.84%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0246"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.40%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.38%  "
